$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row with coal_coke data (normalization/disaggregation addition)
$ws.Range("A6").Value = "coal_coke"
$ws.Range("B6").Value = "coal_coke"

# Move selection to reflect new active cell after edit
$ws.Range("A7").Select()
